# Weekly Fruta/Hortaliza update: a new weekly price observation is
# inserted at row 75 (pushing the existing rows 75-98 down to 76-99).
#
# Insert a new row at position 75; Excel shifts rows 75:98 down to 76:99
# and extends the used range to A1:R99, matching the rest of the sheet's
# row layout/formatting (e.g. the date-formatted style on column D).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(75).Insert()

$ws.Range("A75").Value = 6
$ws.Range("B75").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C75").Value = "Metropolitana"
$ws.Range("D75").Value = 44466
$ws.Range("E75").Value = 13
$ws.Range("F75").Value = 100112001
$ws.Range("G75").Value = "Berenjena"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 410
$ws.Range("K75").Value = 7000
$ws.Range("L75").Value = 8000
$ws.Range("M75").Value = 7439
$ws.Range("N75").Value = "$/caja 50 unidades"
$ws.Range("O75").Value = "Región de Arica y Parinacota"
$ws.Range("P75").Value = 149
$ws.Range("Q75").Value = 50
$ws.Range("R75").Value = "Hortaliza"
